$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 ("tylostyle") values with new measurements
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 460
$ws.Range("D2").Value = 523
$ws.Range("E2").Value = 512
$ws.Range("F2").Value = 501

# Add new row 3 ("raphide") with the values previously held in row 2
$ws.Range("A3").Value = "raphide"
$ws.Range("B3").Value = 118.1
$ws.Range("C3").Value = 126.8
$ws.Range("D3").Value = 135
$ws.Range("E3").Value = 97.7
$ws.Range("F3").Value = 103.7

# Match the author's final selection state
$null = $ws.Range("J12").Select()
